# Update building block types (Imaging assay template)
# - bump template Version 1.0.3 -> 1.0.4
# - rename several "Parameter [...]" building blocks to "Component [...]"
# - rename the Objective term-source/accession columns from OME: to REPR:
# - clear the two "user-specific" term-source values that used to back
#   the (now renamed) microscope / digital camera component columns

$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump version -------------------------------------
$meta = $wb.Worksheets.Item("isa_template")
$meta.Range("B4").Value = "1.0.4"

# --- 3ASY05_Imaging sheet: rename building-block headers -------------------
$data = $wb.Worksheets.Item("3ASY05_Imaging")

$data.Range("E1").Value = "Component [microscope]"
$data.Range("H1").Value = "Component [digital camera]"
$data.Range("K1").Value = "Component [Objective]"
$data.Range("L1").Value = "Term Source REF (REPR:Objective)"
$data.Range("M1").Value = "Term Accession Number (REPR:Objective)"
$data.Range("T1").Value = "Component [Imaging Software Name]"

# --- clear the now-stale "user-specific" values in row 2 -------------------
$data.Range("F2").Value = ""
$data.Range("I2").Value = ""
